# Refresh of the "cryptos" price/volume snapshot (GitHub Actions data pull).
# Updates Price (col D) and Volume(1h) (col E) for most rows, swaps the
# Maker / TheGraph rows (39 <-> 40, new figures too), and replaces the
# FLOKI row (51) with Monero.
#
# Several Price cells hold plain-looking numbers (e.g. "53.00", "0.396")
# but are stored as *text* in the original sheet (trailing zeros / exact
# decimal forms that a numeric cell would not preserve). Assigning those
# directly would make Excel auto-convert them to floating point numbers
# and mangle the formatting (e.g. "53.00" -> 53, "0.0000307" -> 3.07E-05).
# To keep them as text we use Excel's standard leading-apostrophe
# "treat as text" convention (the apostrophe itself is not stored).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.402.16"
$ws.Range("E2").Value = "  -1.89%  "
$ws.Range("D3").Value = "3.487.91"
$ws.Range("E3").Value = "  -2.29%  "
$ws.Range("D5").Value = "`'609.59"
$ws.Range("E5").Value = "  +4.74%  "
$ws.Range("D6").Value = "`'185.97"
$ws.Range("E6").Value = "  -0.10%  "
$ws.Range("E7").Value = "  -0.49%  "
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("E9").Value = "  -4.18%  "
$ws.Range("E10").Value = "  -0.64%  "
$ws.Range("D11").Value = "`'53.00"
$ws.Range("E11").Value = "  -2.73%  "
$ws.Range("D12").Value = "`'0.0000307"
$ws.Range("E12").Value = "  -3.71%  "
$ws.Range("D13").Value = "`'9.52"
$ws.Range("E13").Value = "  +0.61%  "
$ws.Range("D14").Value = "4.043.71"
$ws.Range("E14").Value = "  -2.27%  "
$ws.Range("D15").Value = "`'601.22"
$ws.Range("E15").Value = "  +4.85%  "
$ws.Range("D16").Value = "69.488.72"
$ws.Range("E16").Value = "  -1.85%  "
$ws.Range("D17").Value = "`'18.85"
$ws.Range("E17").Value = "  -2.07%  "
$ws.Range("D18").Value = "`'12.56"
$ws.Range("E18").Value = "  -1.72%  "
$ws.Range("D19").Value = "3.509.32"
$ws.Range("E19").Value = "  -1.71%  "
$ws.Range("E20").Value = "  -0.13%  "
$ws.Range("E21").Value = "  -1.77%  "
$ws.Range("D22").Value = "`'17.14"
$ws.Range("E22").Value = "  -2.56%  "
$ws.Range("D23").Value = "`'105.95"
$ws.Range("E23").Value = "  +12.72%  "
$ws.Range("D24").Value = "`'5.10"
$ws.Range("E24").Value = "  +4.16%  "
$ws.Range("E25").Value = "  +1.65%  "
$ws.Range("E26").Value = "  +3.27%  "
$ws.Range("E27").Value = "  -2.57%  "
$ws.Range("D28").Value = "`'9.67"
$ws.Range("E28").Value = "  +4.98%  "
$ws.Range("D29").Value = "`'33.47"
$ws.Range("E29").Value = "  +3.29%  "
$ws.Range("D30").Value = "`'6.97"
$ws.Range("E30").Value = "  -3.19%  "
$ws.Range("D31").Value = "`'4.14"
$ws.Range("E31").Value = "  +19.02%  "
$ws.Range("D32").Value = "`'12.42"
$ws.Range("E32").Value = "  +1.02%  "
$ws.Range("E33").Value = "  -1.31%  "
$ws.Range("D34").Value = "`'63.35"
$ws.Range("E34").Value = "  +0.42%  "
$ws.Range("E35").Value = "  -6.94%  "
$ws.Range("E36").Value = "  -0.10%  "
$ws.Range("D37").Value = "`'521.59"
$ws.Range("E37").Value = "  -4.97%  "
$ws.Range("D38").Value = "`'3.68"
$ws.Range("E38").Value = "  +7.36%  "
$ws.Range("B39").Value = "TheGraph"
$ws.Range("C39").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D39").Value = "`'0.396"
$ws.Range("E39").Value = "  -4.42%  "
$ws.Range("B40").Value = "Maker"
$ws.Range("C40").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D40").Value = "3.596.88"
$ws.Range("E40").Value = "  +0.38%  "
$ws.Range("D41").Value = "`'36.70"
$ws.Range("E41").Value = "  -3.54%  "
$ws.Range("E42").Value = "  -3.39%  "
$ws.Range("E43").Value = "  -0.98%  "
$ws.Range("E44").Value = "  -0.50%  "
$ws.Range("E45").Value = "  +0.65%  "
$ws.Range("E46").Value = "  +3.14%  "
$ws.Range("D47").Value = "`'3.33"
$ws.Range("E47").Value = "  -4.36%  "
$ws.Range("D48").Value = "`'8.79"
$ws.Range("E49").Value = "  +0.31%  "
$ws.Range("E50").Value = "  -10.19%  "
$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D51").Value = "`'130.45"
$ws.Range("E51").Value = "  -2.93%  "
